# ============================================================================
# Applies the diff: adds w:proofErr grammar-check markers around several
# runs, fixes "К.П. Ситникоа" -> "Ситников К.П.", relocates the _GoBack
# bookmark, retypes "Дата:" / "Например:" run splits, and bumps the cached
# header PAGE field result from 2 to 3.
# ============================================================================

$d = $word.ActiveDocument

# namespace prefix used by every raw-OOXML fragment we hand to InsertXML
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-Pkg([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wns + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
}

# Replace the interior of a paragraph (everything except its trailing pilcrim)
# with a literal run of OOXML, leaving <w:pPr> (and thus all paragraph
# formatting) completely untouched.
function Set-ParaInnerXml($para, [string]$innerXml) {
    $full = $para.Range
    $inner = $d.Range($full.Start, $full.End - 1)
    $inner.InsertXML((New-Pkg ('<w:p>' + $innerXml + '</w:p>')))
}

# ----------------------------------------------------------------------
# Paragraph 6: "Выполнил:_____________"
# ----------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$inner6 = '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Выполнил</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>_</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>____________</w:t></w:r>'
Set-ParaInnerXml $p6 $inner6

# Drop the pre-existing _GoBack bookmark (end of paragraph 7) up front --
# it gets relocated into paragraph 10 below. Hidden bookmarks don't show up
# via enumeration/Count, but are still reachable by name.
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ----------------------------------------------------------------------
# Paragraph 7: "К.П. Ситникоа, группа КЭ-401" -> "Ситников К.П., группа КЭ-401"
# ----------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$inner7 = '<w:r><w:t>Ситников</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> К.П.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, группа </w:t></w:r>' +
    '<w:r><w:t>КЭ</w:t></w:r>' +
    '<w:r><w:t>-</w:t></w:r>' +
    '<w:r><w:t>401</w:t></w:r>'
Set-ParaInnerXml $p7 $inner7

# ----------------------------------------------------------------------
# Paragraph 9: "Руководитель практики:__________"
# ----------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$inner9 = '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Руководитель </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>практики:</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>_</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>_________</w:t></w:r>'
Set-ParaInnerXml $p9 $inner9

# ----------------------------------------------------------------------
# Paragraph 10: "Ст. преподаватель Федянина Р.С." -> "ст. преподаватель Федянина Р.С."
# plus the _GoBack bookmark now lands between "т. преподаватель" and " Федянина"
# ----------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$inner10 = '<w:r><w:t>с</w:t></w:r>' +
    '<w:r><w:t>т. преподаватель</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve"> Федянина</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Р.С.</w:t></w:r>'
Set-ParaInnerXml $p10 $inner10

# ----------------------------------------------------------------------
# Paragraph 12: "Дата:     ________________"
# ----------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$inner12 = '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Дата:   </w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">  ________________</w:t></w:r>'
Set-ParaInnerXml $p12 $inner12

# ----------------------------------------------------------------------
# Paragraph 26: "Например: "
# ----------------------------------------------------------------------
$p26 = $d.Paragraphs(26)
$inner26 = '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:i/><w:color w:val="auto"/><w:szCs w:val="20"/></w:rPr><w:t>Например</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:i/><w:color w:val="auto"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>'
Set-ParaInnerXml $p26 $inner26

# ----------------------------------------------------------------------
# Header: cached PAGE field result "2" -> "3"
# ----------------------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("2", $false, $false, $false, $false, $false, $true, 1, $false, "3", 2) | Out-Null

Write-Host "done"
